{"js": "// Applies the \"added links for capstone\" resume edit:\n//  - Remove the standalone \"GitHub: <link>\" contact line.\n//  - Tighten the Skills & Abilities bullet list text.\n//  - Rewrite three Experience bullet points with more descriptive text.\n//  - Remove one extra blank paragraph before the RESEARCH section.\n//\n// Paragraphs are processed from the bottom of the document upward so that\n// earlier edits (in particular the two paragraph deletions) never shift the\n// index of a paragraph we still need to touch.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map the paragraphs we care about by matching on their current text so the\n// script is resilient to the exact index even if run elsewhere.\nconst items = paragraphs.items;\n\nfunction findIndex(predicate, fromEnd) {\n  if (fromEnd) {\n    for (let i = items.length - 1; i >= 0; i--) {\n      if (predicate(items[i].text)) return i;\n    }\n  } else {\n    for (let i = 0; i < items.length; i++) {\n      if (predicate(items[i].text)) return i;\n    }\n  }\n  return -1;\n}\n\nconst idxGithub = findIndex((t) => t.trim().startsWith(\"GitHub:\"));\nconst idxWebDesign = findIndex((t) => t.indexOf(\"Web Design: HTML and CSS\") !== -1);\nconst idxSkillsList = findIndex((t) => t.trim() === \"JavaScript, Java, SpringBoot, React, MySQL\");\nconst idxBasicPython = findIndex((t) => t.indexOf(\"Basic Python: Introductory Programming\") !== -1);\nconst idxColmar = findIndex((t) => t.indexOf(\"Completed Colmar Academy task\") !== -1);\nconst idxFindYourHat = findIndex((t) => t.indexOf(\"Completed Assessment - Programming a Find Your Hat\") !== -1);\nconst idxCapstone = findIndex((t) => t.indexOf(\"Capstone project - Social Media Website\") !== -1);\n\n// The two identical empty paragraphs sit right after \"Performed diverse\n// administrative tasks\" and before the \"RESEARCH | NUS AVIAN...\" heading;\n// grab the first of the pair.\nconst idxAdmin = findIndex((t) => t.indexOf(\"Performed diverse administrative tasks\") !== -1);\nconst idxEmptyBlank = idxAdmin + 1;\n\n// Sanity: collect every edit as {index, action} then apply from the bottom\n// of the document upward.\nconst edits = [\n  { index: idxEmptyBlank, action: \"delete\" },\n  { index: idxCapstone, action: \"text\", value: \"Completed our capstone project - A social media website for plant lovers using Java Spring Boot for the backend, and JavaScript/HTML & Bootstrap/CSS for the frontend. Information displayed on the front end is stored in a SQL database\" },\n  { index: idxFindYourHat, action: \"text\", value: \"Developed a JavaScript console based game - Find Your Hat\" },\n  { index: idxColmar, action: \"text\", value: \"Developed a mobile reactive webpage based on a specific layout (Colmar Academy webpage)\" },\n  { index: idxBasicPython, action: \"text\", value: \"Basic Python\" },\n  { index: idxSkillsList, action: \"text\", value: \"Java, SpringBoot, MySQL\" },\n  { index: idxWebDesign, action: \"text\", value: \"HTML, CSS, React, and JavaScript\" },\n  { index: idxGithub, action: \"delete\" },\n];\n\nedits.sort((a, b) => b.index - a.index);\n\nfor (const edit of edits) {\n  if (edit.index < 0) continue;\n  const para = paragraphs.items[edit.index];\n  if (edit.action === \"delete\") {\n    para.delete();\n  } else if (edit.action === \"text\") {\n    const range = para.getRange();\n    range.insertText(edit.value, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Applies the \"added links for capstone\" resume edit:\n#  - Remove the standalone \"GitHub: <link>\" contact line.\n#  - Tighten the Skills & Abilities bullet list text.\n#  - Rewrite three Experience bullet points with more descriptive text.\n#  - Remove one extra blank paragraph before the RESEARCH section.\n#\n# Paragraphs are located by matching on their current text (so the script\n# does not depend on hard-coded indices), then every edit is applied from\n# the bottom of the document upward so that the two paragraph deletions\n# never shift the index of a paragraph that still needs to be touched.\n\n$d = $word.ActiveDocument\n\n$idxGithub = -1\n$idxWebDesign = -1\n$idxSkillsList = -1\n$idxBasicPython = -1\n$idxColmar = -1\n$idxFindYourHat = -1\n$idxCapstone = -1\n$idxAdmin = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($idxGithub -lt 0 -and $t.TrimStart().StartsWith(\"GitHub:\")) {\n        $idxGithub = $i\n    }\n    if ($idxWebDesign -lt 0 -and $t.Contains(\"Web Design: HTML and CSS\")) {\n        $idxWebDesign = $i\n    }\n    if ($idxSkillsList -lt 0 -and $t.Trim() -eq \"JavaScript, Java, SpringBoot, React, MySQL\") {\n        $idxSkillsList = $i\n    }\n    if ($idxBasicPython -lt 0 -and $t.Contains(\"Basic Python: Introductory Programming\")) {\n        $idxBasicPython = $i\n    }\n    if ($idxColmar -lt 0 -and $t.Contains(\"Completed Colmar Academy task\")) {\n        $idxColmar = $i\n    }\n    if ($idxFindYourHat -lt 0 -and $t.Contains(\"Completed Assessment - Programming a Find Your Hat\")) {\n        $idxFindYourHat = $i\n    }\n    if ($idxCapstone -lt 0 -and $t.Contains(\"Capstone project - Social Media Website\")) {\n        $idxCapstone = $i\n    }\n    if ($idxAdmin -lt 0 -and $t.Contains(\"Performed diverse administrative tasks\")) {\n        $idxAdmin = $i\n    }\n}\n\n# The two identical empty paragraphs sit right after \"Performed diverse\n# administrative tasks\" and before the \"RESEARCH | NUS AVIAN...\" heading;\n# grab the first of the pair.\n$idxEmptyBlank = $idxAdmin + 1\n\n$edits = @()\n$edits += @{ Index = $idxEmptyBlank; Action = \"delete\" }\n$edits += @{ Index = $idxCapstone; Action = \"text\"; Value = \"Completed our capstone project - A social media website for plant lovers using Java Spring Boot for the backend, and JavaScript/HTML & Bootstrap/CSS for the frontend. Information displayed on the front end is stored in a SQL database\" }\n$edits += @{ Index = $idxFindYourHat; Action = \"text\"; Value = \"Developed a JavaScript console based game - Find Your Hat\" }\n$edits += @{ Index = $idxColmar; Action = \"text\"; Value = \"Developed a mobile reactive webpage based on a specific layout (Colmar Academy webpage)\" }\n$edits += @{ Index = $idxBasicPython; Action = \"text\"; Value = \"Basic Python\" }\n$edits += @{ Index = $idxSkillsList; Action = \"text\"; Value = \"Java, SpringBoot, MySQL\" }\n$edits += @{ Index = $idxWebDesign; Action = \"text\"; Value = \"HTML, CSS, React, and JavaScript\" }\n$edits += @{ Index = $idxGithub; Action = \"delete\" }\n\n$sorted = $edits | Sort-Object -Property Index -Descending\n\nforeach ($edit in $sorted) {\n    if ($edit.Index -lt 1) { continue }\n    $para = $d.Paragraphs.Item($edit.Index)\n    if ($edit.Action -eq \"delete\") {\n        $para.Range.Delete()\n    } else {\n        $para.Range.Text = $edit.Value\n    }\n}\n"}
